# Rotate two blocks of vendor rows in the red-flag report.
# Block 1: rows 2-5  -> new2=old4, new3=old5, new4=old2, new5=old3
# Block 2: rows 8-11 -> new8=old10, new9=old11, new10=old8, new11=old9
# (equivalent to moving the first two rows of each 4-row block to the end)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

function Read-RowData($row) {
    $data = @{}
    foreach ($col in $cols) {
        $data[$col] = $ws.Range($col + $row).Text
    }
    return $data
}

function Write-RowData($row, $data) {
    foreach ($col in $cols) {
        $val = $data[$col]
        $addr = $col + $row
        $current = $ws.Range($addr).Text
        if ($val.Equals($current)) {
            # Nothing actually changes for this cell - leave it untouched
            # so we don't needlessly churn its style/format. NOTE: use the
            # case-sensitive .Equals() - PowerShell's -eq is case-insensitive
            # and would wrongly treat "Rwanda" and "RWANDA" as identical.
            continue
        }
        if ($val.Equals("")) {
            $ws.Range($addr).Value2 = $null
        } else {
            # Force text storage so numeric-looking strings (ids, phone
            # numbers, leading-zero account numbers, ...) survive intact.
            $ws.Range($addr).NumberFormat = "@"
            $ws.Range($addr).Value2 = $val
            $ws.Range($addr).NumberFormat = "General"
        }
    }
}

# Snapshot every source row BEFORE writing anything, so the rotation
# doesn't clobber data it still needs to read later.
$row2 = Read-RowData 2
$row3 = Read-RowData 3
$row4 = Read-RowData 4
$row5 = Read-RowData 5

$row8  = Read-RowData 8
$row9  = Read-RowData 9
$row10 = Read-RowData 10
$row11 = Read-RowData 11

Write-RowData 2 $row4
Write-RowData 3 $row5
Write-RowData 4 $row2
Write-RowData 5 $row3

Write-RowData 8  $row10
Write-RowData 9  $row11
Write-RowData 10 $row8
Write-RowData 11 $row9
